$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$descriptions = @(
    "The Dodge Challenger (2017) is available for `$37,590. It is a reliable choice offering a blend of performance and comfort.",
    "The Ford F-150 XLT (2017) is available for `$30,990. It is a reliable choice offering a blend of performance and comfort.",
    "The Volkswagen Jetta Trendline (2023) is available for `$23,990. It is a reliable choice offering a blend of performance and comfort.",
    "The Honda Accord Touring 2.0 (2020) is available for `$28,990. It is a reliable choice offering a blend of performance and comfort.",
    "The Mazda CX-3 GX (2019) is available for `$21,590. It is a reliable choice offering a blend of performance and comfort.",
    "The Honda Civic Sedan Touring (2022) is available for `$28,590. It is a reliable choice offering a blend of performance and comfort.",
    "The Ford Mustang EcoBoost (2018) is available for `$23,590. It is a reliable choice offering a blend of performance and comfort.",
    "The Ford Mustang GT (2021) is available for `$59,980. It is a reliable choice offering a blend of performance and comfort.",
    "The Hyundai Ioniq 5 (2023) is available for `$42,990. It is a reliable choice offering a blend of performance and comfort.",
    "The Hyundai Elantra (2021) is available for `$21,990. It is a reliable choice offering a blend of performance and comfort.",
    "The Tesla Model 3 (2022) is available for `$38,590. It is a reliable choice offering a blend of performance and comfort.",
    "The Nissan Qashqai S (2021) is available for `$19,990. It is a reliable choice offering a blend of performance and comfort.",
    "The Hyundai Sonata Sport 1.6T (2022) is available for `$24,990. It is a reliable choice offering a blend of performance and comfort.",
    "The Kia Rio EX (2023) is available for `$23,990. It is a reliable choice offering a blend of performance and comfort.",
    "The Tesla Model 3 LR (2023) is available for `$44,990. It is a reliable choice offering a blend of performance and comfort.",
    "The Mazda Mazda 3 Sport GS (2024) is available for `$28,990. It is a reliable choice offering a blend of performance and comfort.",
    "The Ford Mustang V6 Coupe (2016) is available for `$21,790. It is a reliable choice offering a blend of performance and comfort."
)

for ($i = 0; $i -lt $descriptions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $descriptions[$i]
}

$ws.Range("G23").Select() | Out-Null
